$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 11:56"

# Row 4
$ws.Range("B4").Value = 5701285
$ws.Range("C4").Value = 354
$ws.Range("D4").Value = 3063213
$ws.Range("E4").Value = 2461708
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 176364

# Row 6
$ws.Range("B6").Value = 2841400
$ws.Range("C6").Value = 5578
$ws.Range("D6").Value = 2097766
$ws.Range("E6").Value = 689617
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 23
$ws.Range("H6").Value = 54017

# Row 26
$ws.Range("B26").Value = 147211
$ws.Range("C26").Value = 2266
$ws.Range("D26").Value = 100674
$ws.Range("E26").Value = 40119
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 72
$ws.Range("H26").Value = 6418

# Row 32
$ws.Range("A32").Value = "Israel"
$ws.Range("B32").Value = 98550
$ws.Range("C32").Value = 581
$ws.Range("D32").Value = 73848
$ws.Range("E32").Value = 23913
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 8
$ws.Range("H32").Value = 789

# Row 33
$ws.Range("A33").Value = "Ucrania"
$ws.Range("B33").Value = 98537
$ws.Range("C33").Value = 2134
$ws.Range("D33").Value = 50441
$ws.Range("E33").Value = 45912
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 40
$ws.Range("H33").Value = 2184

# Row 38
$ws.Range("A38").Value = "Oman"
$ws.Range("B38").Value = 83769
$ws.Range("C38").Value = 163
$ws.Range("D38").Value = 78386
$ws.Range("E38").Value = 4774
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 6
$ws.Range("H38").Value = 609

# Row 39
$ws.Range("A39").Value = "Panama"
$ws.Range("B39").Value = 83754
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 58274
$ws.Range("E39").Value = 23653
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 1827

# Row 47
$ws.Range("B47").Value = 59378
$ws.Range("C47").Value = 767
$ws.Range("D47").Value = 40481
$ws.Range("E47").Value = 16972
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 12
$ws.Range("H47").Value = 1925

# Row 53
$ws.Range("E53").Value = 3493
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 179

# Row 71
$ws.Range("A71").Value = "Austria"
$ws.Range("B71").Value = 24431
$ws.Range("C71").Value = 347
$ws.Range("D71").Value = 21093
$ws.Range("E71").Value = 2609
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 729

# Row 72
$ws.Range("A72").Value = "Australia"
$ws.Range("B72").Value = 24236
$ws.Range("C72").Value = 243
$ws.Range("D72").Value = 17854
$ws.Range("E72").Value = 5919
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 13
$ws.Range("H72").Value = 463

# Row 84
$ws.Range("B84").Value = 12582
$ws.Range("C84").Value = 36
$ws.Range("D84").Value = 6456
$ws.Range("E84").Value = 5314
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 812

# Row 92
$ws.Range("B92").Value = 9240
$ws.Range("C92").Value = 5
$ws.Range("D92").Value = 8932
$ws.Range("E92").Value = 183
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 125

# Row 98
$ws.Range("A98").Value = "Finlandia"
$ws.Range("B98").Value = 7842
$ws.Range("C98").Value = 37
$ws.Range("D98").Value = 7100
$ws.Range("E98").Value = 408
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 334

# Row 99
$ws.Range("A99").Value = "Albania"
$ws.Range("B99").Value = 7812
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 3928
$ws.Range("E99").Value = 3650
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 234

# Row 128
$ws.Range("A128").Value = "Eslovenia"
$ws.Range("B128").Value = 2536
$ws.Range("C128").Value = 43
$ws.Range("D128").Value = 2079
$ws.Range("E128").Value = 328
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 129

# Row 129
$ws.Range("A129").Value = "Lituania"
$ws.Range("B129").Value = 2528
$ws.Range("C129").Value = 32
$ws.Range("D129").Value = 1747
$ws.Range("E129").Value = 699
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 82

# Row 130
$ws.Range("A130").Value = "Sudan del Sur"
$ws.Range("B130").Value = 2494
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 1290
$ws.Range("E130").Value = 1157
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 47

# Row 152
$ws.Range("B152").Value = 1292
$ws.Range("C152").Value = 7
$ws.Range("D152").Value = 1026
$ws.Range("E152").Value = 211
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 55

# Row 175
$ws.Range("A175").Value = "Papua Nueva Guinea"
$ws.Range("B175").Value = 359
$ws.Range("C175").Value = 12
$ws.Range("D175").Value = 196
$ws.Range("E175").Value = 160
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 3

# Row 176
$ws.Range("A176").Value = "San Martin (Parte Holandesa)"
$ws.Range("B176").Value = 348
$ws.Range("C176").Value = 15
$ws.Range("D176").Value = 133
$ws.Range("E176").Value = 198
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 17

# Row 213
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
